$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing dollar amounts
$ws.Range("B4").Value = 398
$ws.Range("B7").Value = 406

# Add the new transaction note as a new row at the bottom (row 17),
# matching the formatting used by the preceding notes (row 16)
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A17").Value = "10.03.2025 - Out of PO Markkanen karşılığında NSY'ye 1 dolar vermiştir. (398-406)"

# Update selection to match the final state
$ws.Range("B17").Select()
